$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '21.711.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.69%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.538.06'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.27%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("E5").Value = '  +0.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3936'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '

# Row 8
$ws.Range("E8").Value = '  -1.41%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.38'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.66%  '

# Row 10
$ws.Range("E10").Value = '  -1.96%  '

# Row 11
$ws.Range("E11").Value = '  -2.80%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.12%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.779'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.05%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.56%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.635'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.41%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.529.58'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.04%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001095'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.34%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06621'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.24%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.01%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9968'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.10%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.154'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.23%  '

# Row 22
$ws.Range("E22").Value = '  -2.21%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.41%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.366'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.03%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '21.719.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.72%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.384'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.48%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '151.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.94%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.18%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.869'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.707.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.70%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.67%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.090'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.93%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9640'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -9.75%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08102'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.52%  '

# Row 35
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.202'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.06%  '

# Row 36
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.539'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.87%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.498'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.90%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02228'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.50%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05978'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.46%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.34%  '

# Row 41
$ws.Range("E41").Value = '  -1.71%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.181'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.91%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.33%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5824'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.85%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.12'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.51%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.730'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.64%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5590'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.79%  '

# Row 48
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.165'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.56%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.892'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.46%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '115.92'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.35%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06725'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.17%  '
